# Apply edits described by the diff:
# 1. On "Teste1": add row 5 with value "Acura" in A5, and move selection to A5.
# 2. Add new worksheet "Teste2" after "Teste1" with:
#       A1 = "modelo de carro" (header)
#       A2 = "Acura"
#    and selection on A2.
# 3. Set workbook window so the first visible sheet tab is "Teste2" (firstSheet index 1).

$wb = $excel.ActiveWorkbook

# --- Work on Teste1 ---
$ws1 = $wb.Worksheets.Item("Teste1")
$ws1.Range("A5").Value = "Acura"
$ws1.Range("A5").Select()

# --- Add Teste2 right after Teste1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Teste2"
$ws2.Range("A1").Value = "modelo de carro"
$ws2.Range("A2").Value = "Acura"
$ws2.Range("A2").Select()

# Re-activate Teste1 as the selected sheet
$ws1.Activate()
$ws1.Range("A5").Select()

# Scroll the tab bar so Teste2 is the first visible sheet tab
$wb.Windows.Item(1).ScrollWorkbookTabs(1)
